$d = $word.ActiveDocument

$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷7=6, 0", 2)
$d.Content.Find.Execute("48÷6=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷8=2, 6", 2)
$d.Content.Find.Execute("23÷2=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 2)
$d.Content.Find.Execute("24÷9=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=5, 6", 2)
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=3, 3", 2)
$d.Content.Find.Execute("57÷6=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2)
$d.Content.Find.Execute("78÷6=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=30, 0", 2)
$d.Content.Find.Execute("35÷7=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=13, 6", 2)
$d.Content.Find.Execute("64÷7=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 2)
$d.Content.Find.Execute("72÷8=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=3, 4", 2)
$d.Content.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "10÷4=2, 2", 2)
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "68÷7=9, 5", 2)
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=28, 0", 2)
$d.Content.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=46, 1", 2)
$d.Content.Find.Execute("24÷7=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2)
$d.Content.Find.Execute("97÷3=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷4=7, 2", 2)
$d.Content.Find.Execute("52÷4=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷8=10, 2", 2)
$d.Content.Find.Execute("86÷7=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷6=8, 2", 2)
$d.Content.Find.Execute("95÷9=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=13, 1", 2)
$d.Content.Find.Execute("67÷7=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=31, 0", 2)
$d.Content.Find.Execute("17÷3=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=11, 3", 2)
$d.Content.Find.Execute("38÷4=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷2=32, 1", 2)
$d.Content.Find.Execute("21÷2=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=5, 6", 2)
$d.Content.Find.Execute("73÷5=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=5, 4", 2)
$d.Content.Find.Execute("88÷9=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=6, 1", 2)
